# oncogenes now does driver and target
# Adds "Test Done" / "Name" / "comments" data for the "driver gener" (row 37)
# and "target" (row 38) rows, mirroring the pattern used by other rows:
#   C = "yes", D = "oncogenes", E = "db troubles"
# Also updates the sheet's on-screen scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C37").Value = "yes"
$ws.Range("D37").Value = "oncogenes"
$ws.Range("E37").Value = "db troubles"

$ws.Range("C38").Value = "yes"
$ws.Range("D38").Value = "oncogenes"
$ws.Range("E38").Value = "db troubles"

# Move the viewport / selection to match the edited workbook's saved view.
$ws.Range("H34").Select()
